$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel keeps them as literal text (preserving trailing zeros / multi-dot format).
$textFormatRows = @(4,5,6,7,8,9,11,12,13,14,16,17,18,19,20,21,22,24,25,26,27,29,30,32,33,34,35,36,37,38,39,40,41,43,44,45,46,47,48,49,50,51)
foreach ($r in $textFormatRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Update Price column (D)
$ws.Cells.Item(2, 4).Value = "28.086.16"
$ws.Cells.Item(3, 4).Value = "1.818.90"
$ws.Cells.Item(4, 4).Value = "0.9990"
$ws.Cells.Item(5, 4).Value = "310.67"
$ws.Cells.Item(6, 4).Value = "0.9995"
$ws.Cells.Item(7, 4).Value = "0.5008"
$ws.Cells.Item(8, 4).Value = "0.3920"
$ws.Cells.Item(9, 4).Value = "0.1002"
$ws.Cells.Item(11, 4).Value = "40.83"
$ws.Cells.Item(12, 4).Value = "6.421"
$ws.Cells.Item(13, 4).Value = "20.59"
$ws.Cells.Item(14, 4).Value = "0.9985"
$ws.Cells.Item(15, 4).Value = "1.816.44"
$ws.Cells.Item(16, 4).Value = "7.307"
$ws.Cells.Item(17, 4).Value = "0.00001140"
$ws.Cells.Item(18, 4).Value = "92.54"
$ws.Cells.Item(19, 4).Value = "0.06645"
$ws.Cells.Item(20, 4).Value = "0.9985"
$ws.Cells.Item(21, 4).Value = "17.21"
$ws.Cells.Item(22, 4).Value = "5.951"
$ws.Cells.Item(23, 4).Value = "28.143.91"
$ws.Cells.Item(24, 4).Value = "11.12"
$ws.Cells.Item(25, 4).Value = "2.262"
$ws.Cells.Item(26, 4).Value = "158.95"
$ws.Cells.Item(27, 4).Value = "20.78"
$ws.Cells.Item(28, 4).Value = "2.025.98"
$ws.Cells.Item(29, 4).Value = "2.424"
$ws.Cells.Item(30, 4).Value = "127.10"
$ws.Cells.Item(32, 4).Value = "1.040"
$ws.Cells.Item(33, 4).Value = "5.575"
$ws.Cells.Item(34, 4).Value = "3.609"
$ws.Cells.Item(35, 4).Value = "0.06736"
$ws.Cells.Item(36, 4).Value = "0.02343"
$ws.Cells.Item(37, 4).Value = "8.940"
$ws.Cells.Item(38, 4).Value = "0.2139"
$ws.Cells.Item(39, 4).Value = "4.959"
$ws.Cells.Item(40, 4).Value = "11.32"
$ws.Cells.Item(41, 4).Value = "0.6198"
$ws.Cells.Item(43, 4).Value = "0.9987"
$ws.Cells.Item(44, 4).Value = "13.16"
$ws.Cells.Item(45, 4).Value = "0.5920"
$ws.Cells.Item(46, 4).Value = "3.687"
$ws.Cells.Item(47, 4).Value = "1.278"
$ws.Cells.Item(48, 4).Value = "124.44"
$ws.Cells.Item(49, 4).Value = "1.938"
$ws.Cells.Item(50, 4).Value = "1.183"
$ws.Cells.Item(51, 4).Value = "0.06786"

# Update Volume(1h) column (E)
$ws.Cells.Item(2, 5).Value = "  +0.91%  "
$ws.Cells.Item(3, 5).Value = "  +1.61%  "
$ws.Cells.Item(4, 5).Value = "  -0.12%  "
$ws.Cells.Item(5, 5).Value = "  +0.16%  "
$ws.Cells.Item(6, 5).Value = "  -0.05%  "
$ws.Cells.Item(7, 5).Value = "  -2.29%  "
$ws.Cells.Item(8, 5).Value = "  +0.83%  "
$ws.Cells.Item(9, 5).Value = "  +28.22%  "
$ws.Cells.Item(10, 5).Value = "  +1.68%  "
$ws.Cells.Item(11, 5).Value = "  -0.56%  "
$ws.Cells.Item(12, 5).Value = "  +3.19%  "
$ws.Cells.Item(13, 5).Value = "  +1.94%  "
$ws.Cells.Item(14, 5).Value = "  -0.16%  "
$ws.Cells.Item(15, 5).Value = "  +2.30%  "
$ws.Cells.Item(16, 5).Value = "  +1.18%  "
$ws.Cells.Item(17, 5).Value = "  +5.96%  "
$ws.Cells.Item(18, 5).Value = "  +1.10%  "
$ws.Cells.Item(19, 5).Value = "  +1.94%  "
$ws.Cells.Item(20, 5).Value = "  -0.15%  "
$ws.Cells.Item(21, 5).Value = "  +0.87%  "
$ws.Cells.Item(22, 5).Value = "  +0.55%  "
$ws.Cells.Item(23, 5).Value = "  +0.87%  "
$ws.Cells.Item(24, 5).Value = "  +0.89%  "
$ws.Cells.Item(25, 5).Value = "  +1.73%  "
$ws.Cells.Item(26, 5).Value = "  -1.02%  "
$ws.Cells.Item(27, 5).Value = "  +2.52%  "
$ws.Cells.Item(28, 5).Value = "  +1.83%  "
$ws.Cells.Item(29, 5).Value = "  +2.88%  "
$ws.Cells.Item(30, 5).Value = "  +2.50%  "
$ws.Cells.Item(31, 5).Value = "  -0.58%  "
$ws.Cells.Item(32, 5).Value = "  +0.21%  "
$ws.Cells.Item(33, 5).Value = "  +1.49%  "
$ws.Cells.Item(34, 5).Value = "  -0.06%  "
$ws.Cells.Item(35, 5).Value = "  -4.26%  "
$ws.Cells.Item(36, 5).Value = "  +1.67%  "
$ws.Cells.Item(37, 5).Value = "  +2.28%  "
$ws.Cells.Item(38, 5).Value = "  +0.57%  "
$ws.Cells.Item(39, 5).Value = "  -0.67%  "
$ws.Cells.Item(40, 5).Value = "  -1.80%  "
$ws.Cells.Item(41, 5).Value = "  +1.78%  "
$ws.Cells.Item(42, 5).Value = "  +2.66%  "
$ws.Cells.Item(43, 5).Value = "  -0.14%  "
$ws.Cells.Item(44, 5).Value = "  +0.52%  "
$ws.Cells.Item(45, 5).Value = "  +0.38%  "
$ws.Cells.Item(46, 5).Value = "  -0.25%  "
$ws.Cells.Item(47, 5).Value = "  -2.22%  "
$ws.Cells.Item(48, 5).Value = "  +0.05%  "
$ws.Cells.Item(49, 5).Value = "  +1.40%  "
$ws.Cells.Item(50, 5).Value = "  -1.61%  "
$ws.Cells.Item(51, 5).Value = "  -0.47%  "
